$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-48 down to 14-49,
# carrying their values AND formatting with them).
$ws.Rows.Item(13).Insert()

# The freshly inserted row 13 does not inherit the bordered "index column"
# style that column A uses elsewhere, so copy that formatting explicitly
# from the row above before we populate values.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 12 (A8-3) gets updated RPL32/VGR/VG/DCT/DDCT/EXP values.
$row12 = New-Object 'object[,]' 1,12
$row12[0,0] = 'A8-3'
$row12[0,1] = 'A'
$row12[0,2] = 8
$row12[0,3] = 16.88835593149212
$row12[0,4] = 18.1750287128342
$row12[0,5] = 15.78176450473163
$row12[0,6] = -1.106591426760499
$row12[0,7] = 1.286672781342077
$row12[0,8] = 2.65156903710055
$row12[0,9] = 4.42590721110632
$row12[0,10] = 0.1591469002859662
$row12[0,11] = 0.04652315632780536
$ws.Range("B12:M12").Value = $row12

# Row 13 (new) holds the new A8-4 sample.
$ws.Cells.Item(13,1).Value = 11
$row13 = New-Object 'object[,]' 1,12
$row13[0,0] = 'A8-4'
$row13[0,1] = 'A'
$row13[0,2] = 8
$row13[0,3] = 17.60777899023485
$row13[0,4] = 18.59504336307612
$row13[0,5] = 17.62976541333525
$row13[0,6] = 0.02198642310039745
$row13[0,7] = 0.9872643728412669
$row13[0,8] = 3.780146886961447
$row13[0,9] = 4.126498802605511
$row13[0,10] = 0.07278843787437128
$row13[0,11] = 0.05725324192672607
$ws.Range("B13:M13").Value = $row13
